$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data to match the latest scrape.
# For numeric-looking text values in column D, we temporarily force a "Text"
# number format so Excel keeps them as strings (matching the source data,
# which stores these as text, not numbers), then restore the default "Normal"
# style so no stray formatting is left behind on the cell.

$ws.Range("D2").Value = "25.954.76"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.639.91"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06366"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07770"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.281"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("D13").Value = "1.642.47"
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5446"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "0.0₅7817"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").Value = "25.973.20"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "197.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.445"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.948"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.042"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.33%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1173"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.871"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.237"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04997"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.255"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.186"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.538"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.360"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8922"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.584"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.53%  "
$ws.Range("D37").Value = "1.128.18"
$ws.Range("E37").Value = "  -1.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5435"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01557"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("B40").Value = "mCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.550"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.90%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.003"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("E42").Value = "  +10.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.593"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8163"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "1.776.29"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4540"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.000"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.74"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05075"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.003"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.03%  "
